# Traducciones y correcciones menores preinvoice
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edit-fields")

# Fill in rows 14-16 with a new "preinvoice" entity block (same 3-row
# pattern as the other entity/field blocks already present in the sheet:
# type/string, relation/0, trans_name/<field>).

# Row 14: type -> string
$ws.Range("A14").Value = "preinvoice"
$ws.Range("B14").Value = "payment_id"
$ws.Range("C14").Value = "type"
$ws.Range("D14").Value = "string"

# Row 15: relation -> 0
$ws.Range("A15").Value = "preinvoice"
$ws.Range("B15").Value = "payment_id"
$ws.Range("C15").Value = "relation"
$ws.Range("D15").Value = 0

# Row 16: trans_name -> payment_id
$ws.Range("A16").Value = "preinvoice"
$ws.Range("B16").Value = "payment_id"
$ws.Range("C16").Value = "trans_name"
$ws.Range("D16").Value = "payment_id"

# Move the active cell selection to C14, as reflected in the saved view state
$ws.Activate()
$ws.Range("C14").Select()
